# feat: add 2022-Q4 data
#
# The workbook has 4 sheets: "总计" (totals), "2022-Q3", "2022-Q2", "2022-Q1".
# We add a new "2022-Q4" data sheet (copied from the existing "2022-Q3" sheet
# so it inherits the same column layout/formatting), positioned right after
# "总计" and before the other quarters - the existing quarter sheets keep
# their own data untouched, they just shift one tab to the right.
# The "总计" summary sheet gets a new row inserted for 2022-Q4 at the top of
# its data, and the existing rows shift down (their running index in column
# A is renumbered to stay sequential; their own quarter-label / values are
# unchanged).

$wb = $excel.ActiveWorkbook

$totals = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Item("2022-Q3")

# --- 1. Insert the new "2022-Q4" sheet (copy of "2022-Q3" layout), right
#        after "总计" / before the old "2022-Q3" sheet ----------------------
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Fund name (C2) stays the same fund, only the financial figures change.
$q4.Range("D2:G2").NumberFormat = "@"
$q4.Range("D2").Value = "4.76"
$q4.Range("E2").Value = "92.90"
$q4.Range("F2").Value = "3.86"
$q4.Range("G2").Value = "0.1837"
$q4.Range("H2").Value = 7

# --- 2. Update the "总计" sheet: insert a new row for 2022-Q4 ---------------
$totals.Rows.Item(2).Insert()

# Restore the index-column style (border/centering) that Insert() doesn't
# carry onto the brand-new row, and reset the data cells back to the
# workbook's default (unstyled) look used by every other data row.
$totals.Range("A3").Copy()
$totals.Range("A2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$totals.Range("B2:D2").Style = "Normal"

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q4"
$totals.Range("C2").Value = 1
$totals.Range("D2").Value = 0.18

# Renumber the running index in column A for the rows that shifted down.
$totals.Range("A3").Value = 1
$totals.Range("A4").Value = 2
$totals.Range("A5").Value = 3

# Keep "总计" the active/displayed sheet, as it was before the edit.
$totals.Activate()

Write-Host "2022-Q4 sheet + totals row added"
